# "final changes in project"
#
# The AdminUsersPage sheet gets one more row of test data (a "Staff"
# username value) and becomes the sheet that is active/selected when the
# workbook is opened, replacing SubCategoryPage as the selected tab.

$wb = $excel.ActiveWorkbook

# AdminUsersPage is the 4th sheet (rId4 / sheet4.xml) in this workbook.
$ws = $wb.Worksheets.Item("AdminUsersPage")

# Add the new row of data: A3 = "Staff" (grows the used range to A1:B3
# and adds a new shared string).
$ws.Range("A3").Value = "Staff"

# Make AdminUsersPage the active/selected sheet (was SubCategoryPage),
# with A3 - the newly entered cell - selected.
$ws.Activate()
$ws.Range("A3").Select()
